$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44656
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112021
$ws.Cells.Item($row, 7).Value = "Ají"
$ws.Cells.Item($row, 8).Value = "Americana (o)"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 30
$ws.Cells.Item($row, 11).Value = 18000
$ws.Cells.Item($row, 12).Value = 20000
$ws.Cells.Item($row, 13).Value = 19333
$ws.Cells.Item($row, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 773
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
